$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) FirstParagraph: extend the last run's text with new sentences
#    and fix "haft" -> "varit med om".
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    ", en upplevelse han haft många gånger på jobbet.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    ", en upplevelse han varit med om många gånger på jobbet. Fordon körde förbi på vägarna, regn föll från himlen med små droppar som täckte gatorna. Adam, med sina händer i fickorna, funderade på hur många spår i fall som möjligen försvunnit under regnets gång.",
    2) | Out-Null

# ---------------------------------------------------------------
# 2) Remove the old "Fordon körde förbi..." paragraph: its (edited)
#    text now lives inside the FirstParagraph above.
# ---------------------------------------------------------------
$oldPara2 = $d.Paragraphs(3)
$oldPara2.Range.Delete() | Out-Null

# ---------------------------------------------------------------
# 3) "Eventuellt kom Adam fram..." paragraph is unchanged.
# ---------------------------------------------------------------

# ---------------------------------------------------------------
# 4) Door paragraph: fix "instrution" -> "instruktion" and splice in
#    a new sentence about the door/lock before the final sentence.
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    "på, en instrution för städerskana på hotellet. Efter en kort stund så öppnade inspektören dörren och tog försiktiga steg in.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "på, en instruktion för städerskana på hotellet. Locket på dörren verkade se bra ut, ingen försökte ta sig in utan nyckel. Efter en kort stund så öppnade inspektören dörren och tog försiktiga steg in.",
    2) | Out-Null

# ---------------------------------------------------------------
# 5) Replace the "tejp på golvet" paragraph's content with the new
#    "room contents" paragraph text.
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    "Det fanns inte mycket och se i det rätt lilla rummet, vitt tejp som markerade vart kroppen hittades fanns synligt på golvet, som om den skildiga inte hade tid och gömma kroppen. En fläkt satt på uppe på taket stilla.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Mycket fanns det inte i rummet på första blick, det främsta innehållet av rummet log på golvet. Kuddar, lampor, täcken, det såg verkligen ut som om någon haft strid här inne. Lampan glänsde på taket, där det också satt en fläkt, passande för vädret.",
    2) | Out-Null

# ---------------------------------------------------------------
# 6) New paragraph: tape / body position on the floor.
# ---------------------------------------------------------------
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter() | Out-Null
$p6 = $d.Paragraphs($d.Paragraphs.Count)
$p6.Style = "BodyText"
$p6.Range.Text = "På golvet fanns också bitar av vit tejp, fast sätta för att visa vart kroppen hittades, mycket kunde inte synnas, men positionen av kroppen — i mitten av rummet — verkades som en viktig detalj."

# ---------------------------------------------------------------
# 7) New paragraph: windows in the room.
# ---------------------------------------------------------------
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter() | Out-Null
$p7 = $d.Paragraphs($d.Paragraphs.Count)
$p7.Style = "BodyText"
$p7.Range.Text = "Dem få fönstrerna i rummet var alla stängda, och verkade inte kunna öppnas heller, så ingen tog sig in därifrån heller."

# ---------------------------------------------------------------
# 8) New paragraph: the report about "Pearl White" (plain text for
#    now -- the italics on the name are applied at the very end, see
#    step 10, to avoid the new-paragraph formatting leaking forward
#    into paragraph 9).
# ---------------------------------------------------------------
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter() | Out-Null
$p8 = $d.Paragraphs($d.Paragraphs.Count)
$p8.Style = "BodyText"
$r8 = $p8.Range
$r8.Text = "Adam minndes några av de viktigare punkterna kring fallet, ändå tog han fram rapporten kring vad som hittades. Dödsoffret,"
$r8.InsertAfter(" ")
$r8.InsertAfter("“")
$r8.InsertAfter("Pearl White")
$r8.InsertAfter("”")
$r8.InsertAfter(", var en kvinna på besök här från england. Rapporten beskrev att kvinnan avled kring klockan 3 igår på eftermiddagen, inga tydliga skador hittades på kroppen. Mer sidor i rapporten beskrev möjliga teorier kring dödsskälet, det mest möjliga var att kvinnan hade kvävdes till döds.")

# ---------------------------------------------------------------
# 9) New paragraph: Adam checks his watch and heads back.
# ---------------------------------------------------------------
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter() | Out-Null
$p9 = $d.Paragraphs($d.Paragraphs.Count)
$p9.Style = "BodyText"
$p9.Range.Text = "Efter att ha stängt blocket med rapportets innehåll kollade Adam på sin klocka, han var redan sen till att hålla förhör med de misstänkte, men han hade knappt hunnit titta genom brottsplatsen själv. Han bestämde sig eventuellt och acceptera det han hade, och återvända till polisstationen."

# ---------------------------------------------------------------
# 10) Now that no more new paragraphs will be created, italicize the
#     victim's name "Pearl White" inside paragraph 8.
# ---------------------------------------------------------------
$nameRange = $d.Content
$nameRange.Find.Execute("Pearl White", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$nameRange.Font.Italic = $true

Write-Output "done"
foreach ($p in $d.Paragraphs) {
    Write-Output "----"
    Write-Output $p.Range.Text
}
